# Add the new "managerId" column (K) to the questions template sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell K1 — same bold/centered header style as the rest of row 1.
$ws.Range("K1").Value = "managerId"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats

# Data cells K2:K4.
$ws.Range("K2").Value = 123
$ws.Range("K3").Value = 2
$ws.Range("K4").Value = 3

# Clear clipboard marquee and match the saved selection/active cell.
$excel.CutCopyMode = $false
$ws.Range("K4").Select()
